# Refactor projected production to use ghg_s1s2 values and benchmarks:
# zero out the first projection-year (column C) production growth rates
# on the "projected_production" sheet for every sector/region row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("projected_production")

$ws.Range("C2:C7").Value = 0
